# Updates the crypto price/volume table (and two name-swap pairs in the
# rank-9/10 and rank-15..22 ranges) to the new "Updated symbol list" snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume 1h) hold numeric-looking text ("321.70", "-2.29%").
# The source file stores them as plain text, so we write them with a leading
# apostrophe to force text (not a Number/Percentage), then restore the cell's
# original "Normal" style so no extra number formatting is picked up.
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

# row 2
Set-TextValue "D2" "319.93"
Set-TextValue "E2" "-2.85%"

# row 3
Set-TextValue "D3" "42.66"
Set-TextValue "E3" "-4.65%"

# row 4
Set-TextValue "D4" "5.192"
Set-TextValue "E4" "-5.41%"

# row 5
Set-TextValue "D5" "0.08186"
Set-TextValue "E5" "-2.58%"

# row 6
Set-TextValue "D6" "4.365"
Set-TextValue "E6" "-2.09%"

# row 7
Set-TextValue "D7" "1.770"
Set-TextValue "E7" "-12.17%"

# row 8
Set-TextValue "D8" "0.9499"
Set-TextValue "E8" "-3.05%"

# row 9
Set-TextValue "E9" "0.90%"

# row 10
Set-TextValue "D10" "0.1877"
Set-TextValue "E10" "-1.90%"

# row 11 -> BitrueCoin
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.04695"
Set-TextValue "E11" "-0.10%"

# row 12 -> MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09355"
Set-TextValue "E12" "-3.44%"

# row 13
Set-TextValue "D13" "7.431"
Set-TextValue "E13" "-21.75%"

# row 14
Set-TextValue "D14" "0.1059"
Set-TextValue "E14" "-0.07%"

# row 15
Set-TextValue "D15" "0.001284"
Set-TextValue "E15" "-0.13%"

# row 16
Set-TextValue "D16" "0.005697"
Set-TextValue "E16" "-4.81%"

# row 17 -> HotbitToken
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D17" "0.004294"
Set-TextValue "E17" "-3.05%"

# row 18 -> LEO
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.353"
Set-TextValue "E18" "-1.02%"

# row 19 -> BTSEToken
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D19" "2.532"
Set-TextValue "E19" "-0.19%"

# row 20 -> BitpandaEcosystemToken
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D20" "0.3362"
Set-TextValue "E20" "0.34%"

# row 21 -> ProBitToken
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1337"
Set-TextValue "E21" "-1.39%"

# row 22 -> ZBToken
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D22" "0.2557"
Set-TextValue "E22" "0.03%"

# row 23 -> CoinExToken
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D23" "0.04185"
Set-TextValue "E23" "0.05%"

# row 24 -> BitKan
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D24" "0.001242"
Set-TextValue "E24" "-4.70%"

# row 25
Set-TextValue "D25" "0.0001224"
Set-TextValue "E25" "-6.37%"

# row 26
Set-TextValue "D26" "0.0002988"
Set-TextValue "E26" "0.03%"

# row 38
Set-TextValue "D38" "0.02637"
Set-TextValue "E38" "-2.83%"

# row 39
Set-TextValue "D39" "0.05625"
Set-TextValue "E39" "-0.11%"

# row 40
Set-TextValue "D40" "0.008173"
Set-TextValue "E40" "3.72%"

# row 41
Set-TextValue "D41" "0.1401"
Set-TextValue "E41" "-1.73%"

# row 42
Set-TextValue "D42" "0.006563"
Set-TextValue "E42" "-11.19%"

# row 43
Set-TextValue "D43" "0.002118"
Set-TextValue "E43" "-0.41%"

# row 44
Set-TextValue "D44" "0.007617"
Set-TextValue "E44" "-12.03%"

# row 45
Set-TextValue "D45" "0.3475"
Set-TextValue "E45" "2.66%"

# row 46
Set-TextValue "D46" "0.00006792"
Set-TextValue "E46" "-1.45%"

# row 47
Set-TextValue "E47" "-0.24%"

# row 48
Set-TextValue "D48" "0.003348"
Set-TextValue "E48" "-4.21%"

# row 49
Set-TextValue "D49" "0.004113"
Set-TextValue "E49" "15.88%"

# row 50
Set-TextValue "D50" "0.00002107"
Set-TextValue "E50" "-0.24%"

# row 51
Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "-0.24%"
